$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update time values: seconds -> milliseconds (row 3 = time_original, row 10 = time_prune_s1, row 17 = time_prune_s2)
$ws.Range("C3").Value = 0.3126965
$ws.Range("D3").Value = 3.195789401041667
$ws.Range("E3").Value = 3.619022005208334
$ws.Range("F3").Value = 4.967390161458333
$ws.Range("G3").Value = 7.415865105902777
$ws.Range("H3").Value = 7.655040432291666
$ws.Range("I3").Value = 8.828057548611111
$ws.Range("J3").Value = 11.02605521180556
$ws.Range("K3").Value = 13.35302809027778
$ws.Range("L3").Value = 13.49793984895833
$ws.Range("M3").Value = 14.66949730555556
$ws.Range("N3").Value = 16.87689825694444
$ws.Range("O3").Value = 19.20576161111111
$ws.Range("P3").Value = 19.28444795138888
$ws.Range("Q3").Value = 21.15322092361111
$ws.Range("R3").Value = 23.00321997222222
$ws.Range("S3").Value = 24.80391318402777
$ws.Range("T3").Value = 24.824139

$ws.Range("C10").Value = 0.1613542378472222
$ws.Range("D10").Value = 0.7036061701388888
$ws.Range("E10").Value = 0.8464881944444443
$ws.Range("F10").Value = 1.129266887152778
$ws.Range("G10").Value = 1.618272729166667
$ws.Range("H10").Value = 1.687166574652778
$ws.Range("I10").Value = 1.867558008680556
$ws.Range("J10").Value = 2.122761454861111
$ws.Range("K10").Value = 2.329663388888889
$ws.Range("L10").Value = 2.353591822916667
$ws.Range("M10").Value = 2.439550385416667
$ws.Range("N10").Value = 2.553520779513889
$ws.Range("O10").Value = 2.649179206597222
$ws.Range("P10").Value = 2.659316967013889
$ws.Range("Q10").Value = 2.737239118055555
$ws.Range("R10").Value = 2.813186053819444
$ws.Range("S10").Value = 2.887139371527779
$ws.Range("T10").Value = 2.889305788194445

$ws.Range("C17").Value = 0.1577050590277778
$ws.Range("D17").Value = 0.6229637378472223
$ws.Range("E17").Value = 0.7620447170138889
$ws.Range("F17").Value = 0.9985711059027776
$ws.Range("G17").Value = 1.429373045138889
$ws.Range("H17").Value = 1.488095822916667
$ws.Range("I17").Value = 1.640996751736111
$ws.Range("J17").Value = 1.868069293402777
$ws.Range("K17").Value = 2.054396899305555
$ws.Range("L17").Value = 2.075141029513889
$ws.Range("M17").Value = 2.156682425347222
$ws.Range("N17").Value = 2.263981713541667
$ws.Range("O17").Value = 2.3581944375
$ws.Range("P17").Value = 2.367740720486111
$ws.Range("Q17").Value = 2.440125642361112
$ws.Range("R17").Value = 2.5139573125
$ws.Range("S17").Value = 2.580146041666667
$ws.Range("T17").Value = 2.582103788194444

# Update chart axis title: "time elapsed (s)" -> "time elapsed (ms)"
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2, 1)  # xlValue=2, xlPrimary=1
$valueAxis.AxisTitle.Text = "time elapsed (ms)"
